# Update "想去人数" (number of people wanting to go) values in F column
# for both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value
$updates = @{
    "F3"  = 1457
    "F4"  = 952
    "F6"  = 2144
    "F7"  = 35
    "F8"  = 1301
    "F10" = 123
    "F11" = 37
    "F12" = 310
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
